# "Add files via upload" - adds new HTML/CSS rows (font-weight, font-style,
# line-height, text-decoration) plus explanatory comments, re-uploaded by
# user "adm".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 4 new rows right above the existing "HTML Formatting" section
# (currently row 26) so the new attribute rows land at 26-29 and everything
# below shifts down by 4. ---
$ws.Rows.Item(26).Resize(4).Insert()

# Copy the formatting of the row just above (the last "text-align" row,
# style index 2 / thin border) into the freshly inserted blank rows so they
# match the sheet's normal row styling instead of being left unstyled.
$ws.Range("A25:B25").Copy()
$ws.Range("A26:B29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the new attribute / description pairs ---
$ws.Range("A26").Value = "font-weight"
$ws.Range("B26").Value = "Defines for the bold"

$ws.Range("A27").Value = "font-style"
$ws.Range("B27").Value = "for the italic"

$ws.Range("A28").Value = "line-height"
$ws.Range("B28").Value = "distance between the lines"

$ws.Range("A29").Value = "text-decoration"
$ws.Range("B29").Value = "for defining underline"

# --- New explanatory comments from "adm" on the line-height and
# text-decoration rows ---
$c1 = $ws.Range("B28").AddComment()
$c1.Text("adm:" + "`n" + "Can be defined as below" + "`n" + "line-height: 15px" + "`n" + "line-height : 1.5em")
$c1.Visible = $false

$c2 = $ws.Range("B29").AddComment()
$c2.Text("adm:" + "`n" + "text-decoration: underline;" + "`n" + "text-decoration: none to remove the underline from link")
$c2.Visible = $false

# --- Restore the selection/active cell the author left the sheet on ---
$ws.Range("B38").Select()
